# EI Variable Installments T2 scenarios
#
# Adds a new "waittopageload1" / 2000 step as row 8 on the
# "Edit Repayment Schedule" sheet (pushing the existing rows 8-14 down to
# 9-15), and leaves that sheet as the active sheet/selection so the
# workbook reopens on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row above the current row 8 ("clickonsubmit"/"Submit"),
# shifting rows 8:14 down to 9:15.
$ws.Range("A8").EntireRow.Insert() | Out-Null

# Populate the newly inserted row.
$ws.Range("A8").Value = "waittopageload1"
$ws.Range("B8").Value = 2000

# Match the formatting used by the other "click wait" numeric cell (B3)
# instead of the default row styling the insert leaves behind.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

# Make "Edit Repayment Schedule" the active/selected sheet with A8:B8
# selected (this also clears tabSelected on whichever sheet was active
# before, e.g. "NewLoanInput").
$ws.Activate() | Out-Null
$ws.Range("A8:B8").Select() | Out-Null
